$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 128-209: each (Primera, Segunda) pair shifts down by one pair,
# with two brand-new rows (128-129) inserted at the top of this block; the old last pair
# (previously rows 208-209) moves down to become the new final rows 210-211.

$ws.Range("D128").Value = 44438 ; $ws.Range("J128").Value = 2000 ; $ws.Range("K128").Value = 7500 ; $ws.Range("L128").Value = 8000 ; $ws.Range("M128").Value = 7750 ; $ws.Range("P128").Value = 1292
$ws.Range("D129").Value = 44438 ; $ws.Range("J129").Value = 1400 ; $ws.Range("K129").Value = 6000 ; $ws.Range("L129").Value = 7000 ; $ws.Range("M129").Value = 6500 ; $ws.Range("P129").Value = 1083
$ws.Range("D130").Value = 44260 ; $ws.Range("J130").Value = 3200 ; $ws.Range("K130").Value = 6500 ; $ws.Range("L130").Value = 7000 ; $ws.Range("M130").Value = 6750 ; $ws.Range("P130").Value = 1125
$ws.Range("D131").Value = 44260 ; $ws.Range("J131").Value = 1680 ; $ws.Range("K131").Value = 5500 ; $ws.Range("L131").Value = 6000 ; $ws.Range("M131").Value = 5750 ; $ws.Range("P131").Value = 958
$ws.Range("D132").Value = 44209 ; $ws.Range("J132").Value = 3000 ; $ws.Range("K132").Value = 7800 ; $ws.Range("L132").Value = 8000 ; $ws.Range("M132").Value = 7900 ; $ws.Range("P132").Value = 1317
$ws.Range("D133").Value = 44209 ; $ws.Range("J133").Value = 1600 ; $ws.Range("K133").Value = 6800 ; $ws.Range("L133").Value = 7000 ; $ws.Range("M133").Value = 6900 ; $ws.Range("P133").Value = 1150
$ws.Range("D134").Value = 44330 ; $ws.Range("J134").Value = 3240 ; $ws.Range("K134").Value = 9000 ; $ws.Range("L134").Value = 10000 ; $ws.Range("M134").Value = 9500 ; $ws.Range("P134").Value = 1583
$ws.Range("D135").Value = 44330 ; $ws.Range("J135").Value = 1640 ; $ws.Range("K135").Value = 7500 ; $ws.Range("L135").Value = 8000 ; $ws.Range("M135").Value = 7750 ; $ws.Range("P135").Value = 1292
$ws.Range("D136").Value = 44391 ; $ws.Range("J136").Value = 3360 ; $ws.Range("K136").Value = 8000 ; $ws.Range("L136").Value = 9000 ; $ws.Range("M136").Value = 8500 ; $ws.Range("P136").Value = 1417
$ws.Range("D137").Value = 44391 ; $ws.Range("K137").Value = 6000 ; $ws.Range("M137").Value = 6500 ; $ws.Range("P137").Value = 1083
$ws.Range("D138").Value = 44193 ; $ws.Range("J138").Value = 2600 ; $ws.Range("K138").Value = 7800 ; $ws.Range("L138").Value = 8000 ; $ws.Range("M138").Value = 7900 ; $ws.Range("P138").Value = 1317
$ws.Range("D139").Value = 44193 ; $ws.Range("K139").Value = 6800 ; $ws.Range("M139").Value = 6900 ; $ws.Range("P139").Value = 1150
$ws.Range("D140").Value = 44351 ; $ws.Range("J140").Value = 3280 ; $ws.Range("K140").Value = 8000 ; $ws.Range("M140").Value = 8500 ; $ws.Range("P140").Value = 1417
$ws.Range("D141").Value = 44351
$ws.Range("D142").Value = 44358 ; $ws.Range("K142").Value = 8500 ; $ws.Range("M142").Value = 8750 ; $ws.Range("P142").Value = 1458
$ws.Range("D143").Value = 44358 ; $ws.Range("J143").Value = 1600 ; $ws.Range("K143").Value = 6500 ; $ws.Range("M143").Value = 6750 ; $ws.Range("P143").Value = 1125
$ws.Range("D144").Value = 44389 ; $ws.Range("K144").Value = 8000 ; $ws.Range("L144").Value = 9000 ; $ws.Range("M144").Value = 8500 ; $ws.Range("P144").Value = 1417
$ws.Range("D145").Value = 44389 ; $ws.Range("J145").Value = 1400 ; $ws.Range("L145").Value = 7000 ; $ws.Range("M145").Value = 6500 ; $ws.Range("P145").Value = 1083
$ws.Range("D146").Value = 44251 ; $ws.Range("J146").Value = 3200 ; $ws.Range("K146").Value = 7000 ; $ws.Range("L146").Value = 7500 ; $ws.Range("M146").Value = 7250 ; $ws.Range("P146").Value = 1208
$ws.Range("D147").Value = 44251 ; $ws.Range("J147").Value = 1600 ; $ws.Range("K147").Value = 6000 ; $ws.Range("L147").Value = 6500 ; $ws.Range("M147").Value = 6250 ; $ws.Range("P147").Value = 1042
$ws.Range("D148").Value = 44305 ; $ws.Range("J148").Value = 3000 ; $ws.Range("K148").Value = 7500 ; $ws.Range("L148").Value = 8000 ; $ws.Range("M148").Value = 7750 ; $ws.Range("P148").Value = 1292
$ws.Range("D149").Value = 44305 ; $ws.Range("J149").Value = 1480
$ws.Range("D150").Value = 44417 ; $ws.Range("J150").Value = 3200
$ws.Range("D151").Value = 44417 ; $ws.Range("K151").Value = 6500 ; $ws.Range("M151").Value = 6750 ; $ws.Range("P151").Value = 1125
$ws.Range("D152").Value = 44419 ; $ws.Range("J152").Value = 3400 ; $ws.Range("K152").Value = 8000 ; $ws.Range("L152").Value = 9000 ; $ws.Range("M152").Value = 8500 ; $ws.Range("P152").Value = 1417
$ws.Range("D153").Value = 44419 ; $ws.Range("J153").Value = 1600 ; $ws.Range("K153").Value = 6000 ; $ws.Range("M153").Value = 6500 ; $ws.Range("P153").Value = 1083
$ws.Range("D154").Value = 44202 ; $ws.Range("J154").Value = 2800 ; $ws.Range("K154").Value = 7800 ; $ws.Range("M154").Value = 7900 ; $ws.Range("P154").Value = 1317
$ws.Range("D155").Value = 44202 ; $ws.Range("J155").Value = 1560 ; $ws.Range("K155").Value = 6800 ; $ws.Range("M155").Value = 6900 ; $ws.Range("P155").Value = 1150
$ws.Range("D156").Value = 44307 ; $ws.Range("J156").Value = 3200
$ws.Range("D157").Value = 44307 ; $ws.Range("J157").Value = 1680
$ws.Range("D158").Value = 44195 ; $ws.Range("J158").Value = 2800 ; $ws.Range("K158").Value = 7500 ; $ws.Range("L158").Value = 8000 ; $ws.Range("M158").Value = 7750 ; $ws.Range("P158").Value = 1292
$ws.Range("D159").Value = 44195 ; $ws.Range("J159").Value = 1560 ; $ws.Range("K159").Value = 6500 ; $ws.Range("L159").Value = 7000 ; $ws.Range("M159").Value = 6750 ; $ws.Range("P159").Value = 1125
$ws.Range("D160").Value = 44265 ; $ws.Range("K160").Value = 6500 ; $ws.Range("L160").Value = 7000 ; $ws.Range("M160").Value = 6750 ; $ws.Range("P160").Value = 1125
$ws.Range("D161").Value = 44265 ; $ws.Range("J161").Value = 1720 ; $ws.Range("K161").Value = 5500 ; $ws.Range("L161").Value = 6000 ; $ws.Range("M161").Value = 5750 ; $ws.Range("P161").Value = 958
$ws.Range("D162").Value = 44333 ; $ws.Range("J162").Value = 3200 ; $ws.Range("K162").Value = 9000 ; $ws.Range("L162").Value = 9500 ; $ws.Range("M162").Value = 9250 ; $ws.Range("P162").Value = 1542
$ws.Range("D163").Value = 44333 ; $ws.Range("J163").Value = 1460 ; $ws.Range("K163").Value = 7500 ; $ws.Range("L163").Value = 8000 ; $ws.Range("M163").Value = 7750 ; $ws.Range("P163").Value = 1292
$ws.Range("D164").Value = 44277 ; $ws.Range("J164").Value = 2700 ; $ws.Range("K164").Value = 7000 ; $ws.Range("L164").Value = 8000 ; $ws.Range("M164").Value = 7500 ; $ws.Range("P164").Value = 1250
$ws.Range("D165").Value = 44277 ; $ws.Range("K165").Value = 5000 ; $ws.Range("L165").Value = 6000 ; $ws.Range("M165").Value = 5500 ; $ws.Range("P165").Value = 917
$ws.Range("D166").Value = 44433 ; $ws.Range("J166").Value = 2000 ; $ws.Range("K166").Value = 8000 ; $ws.Range("L166").Value = 9000 ; $ws.Range("M166").Value = 8500 ; $ws.Range("P166").Value = 1417
$ws.Range("D167").Value = 44433 ; $ws.Range("J167").Value = 1400 ; $ws.Range("K167").Value = 6500 ; $ws.Range("M167").Value = 6750 ; $ws.Range("P167").Value = 1125
$ws.Range("D168").Value = 44309 ; $ws.Range("J168").Value = 3200 ; $ws.Range("K168").Value = 7800 ; $ws.Range("L168").Value = 8000 ; $ws.Range("M168").Value = 7900 ; $ws.Range("P168").Value = 1317
$ws.Range("D169").Value = 44309 ; $ws.Range("J169").Value = 1660 ; $ws.Range("K169").Value = 6800 ; $ws.Range("L169").Value = 7000 ; $ws.Range("M169").Value = 6900 ; $ws.Range("P169").Value = 1150
$ws.Range("D170").Value = 44344 ; $ws.Range("J170").Value = 3280 ; $ws.Range("K170").Value = 8500 ; $ws.Range("L170").Value = 9000 ; $ws.Range("M170").Value = 8750 ; $ws.Range("P170").Value = 1458
$ws.Range("D171").Value = 44344 ; $ws.Range("J171").Value = 1600 ; $ws.Range("K171").Value = 7000 ; $ws.Range("L171").Value = 7500 ; $ws.Range("M171").Value = 7250 ; $ws.Range("P171").Value = 1208
$ws.Range("D172").Value = 44319 ; $ws.Range("J172").Value = 3000
$ws.Range("D173").Value = 44319 ; $ws.Range("J173").Value = 1480
$ws.Range("D174").Value = 44316 ; $ws.Range("J174").Value = 3400 ; $ws.Range("K174").Value = 7500 ; $ws.Range("L174").Value = 8000 ; $ws.Range("M174").Value = 7750 ; $ws.Range("P174").Value = 1292
$ws.Range("D175").Value = 44316 ; $ws.Range("K175").Value = 6500 ; $ws.Range("L175").Value = 7000 ; $ws.Range("M175").Value = 6750 ; $ws.Range("P175").Value = 1125
$ws.Range("D176").Value = 44253 ; $ws.Range("J176").Value = 3300
$ws.Range("D177").Value = 44253 ; $ws.Range("J177").Value = 1680 ; $ws.Range("K177").Value = 5500 ; $ws.Range("M177").Value = 5750 ; $ws.Range("P177").Value = 958
$ws.Range("D178").Value = 44281 ; $ws.Range("J178").Value = 3200 ; $ws.Range("K178").Value = 7000 ; $ws.Range("L178").Value = 7500 ; $ws.Range("M178").Value = 7250 ; $ws.Range("P178").Value = 1208
$ws.Range("D179").Value = 44281 ; $ws.Range("J179").Value = 1660 ; $ws.Range("K179").Value = 5000 ; $ws.Range("L179").Value = 6000 ; $ws.Range("M179").Value = 5500 ; $ws.Range("P179").Value = 917
$ws.Range("D180").Value = 44160 ; $ws.Range("J180").Value = 2800 ; $ws.Range("K180").Value = 7500 ; $ws.Range("M180").Value = 7750 ; $ws.Range("P180").Value = 1292
$ws.Range("D181").Value = 44160 ; $ws.Range("J181").Value = 1600 ; $ws.Range("K181").Value = 6500 ; $ws.Range("M181").Value = 6750 ; $ws.Range("P181").Value = 1125
$ws.Range("D182").Value = 44186 ; $ws.Range("J182").Value = 2700
$ws.Range("D183").Value = 44186 ; $ws.Range("J183").Value = 1540 ; $ws.Range("K183").Value = 6800 ; $ws.Range("L183").Value = 7000 ; $ws.Range("M183").Value = 6900 ; $ws.Range("P183").Value = 1150
$ws.Range("D184").Value = 44211 ; $ws.Range("J184").Value = 3000 ; $ws.Range("K184").Value = 7800 ; $ws.Range("M184").Value = 7900 ; $ws.Range("P184").Value = 1317
$ws.Range("D185").Value = 44211 ; $ws.Range("J185").Value = 1600 ; $ws.Range("K185").Value = 7000 ; $ws.Range("L185").Value = 7500 ; $ws.Range("M185").Value = 7250 ; $ws.Range("P185").Value = 1208
$ws.Range("D186").Value = 44272 ; $ws.Range("J186").Value = 3300 ; $ws.Range("K186").Value = 7000 ; $ws.Range("M186").Value = 7500 ; $ws.Range("P186").Value = 1250
$ws.Range("D187").Value = 44272 ; $ws.Range("J187").Value = 1680 ; $ws.Range("K187").Value = 5000 ; $ws.Range("M187").Value = 5500 ; $ws.Range("P187").Value = 917
$ws.Range("D188").Value = 44370 ; $ws.Range("J188").Value = 3400 ; $ws.Range("K188").Value = 7500 ; $ws.Range("L188").Value = 8000 ; $ws.Range("M188").Value = 7750 ; $ws.Range("P188").Value = 1292
$ws.Range("D189").Value = 44370 ; $ws.Range("J189").Value = 1640 ; $ws.Range("K189").Value = 5500 ; $ws.Range("L189").Value = 6000 ; $ws.Range("M189").Value = 5750 ; $ws.Range("P189").Value = 958
$ws.Range("D190").Value = 44326 ; $ws.Range("J190").Value = 3100 ; $ws.Range("K190").Value = 9000 ; $ws.Range("L190").Value = 10000 ; $ws.Range("M190").Value = 9500 ; $ws.Range("P190").Value = 1583
$ws.Range("D191").Value = 44326 ; $ws.Range("J191").Value = 1400
$ws.Range("D192").Value = 44230 ; $ws.Range("J192").Value = 3000 ; $ws.Range("K192").Value = 8500 ; $ws.Range("L192").Value = 9000 ; $ws.Range("M192").Value = 8750 ; $ws.Range("P192").Value = 1458
$ws.Range("D193").Value = 44230 ; $ws.Range("J193").Value = 1800 ; $ws.Range("K193").Value = 7500 ; $ws.Range("L193").Value = 8000 ; $ws.Range("M193").Value = 7750 ; $ws.Range("P193").Value = 1292
$ws.Range("D194").Value = 44382 ; $ws.Range("J194").Value = 3200 ; $ws.Range("K194").Value = 7500 ; $ws.Range("L194").Value = 8000 ; $ws.Range("M194").Value = 7750 ; $ws.Range("P194").Value = 1292
$ws.Range("D195").Value = 44382 ; $ws.Range("J195").Value = 1520 ; $ws.Range("K195").Value = 5500 ; $ws.Range("L195").Value = 6000 ; $ws.Range("M195").Value = 5750 ; $ws.Range("P195").Value = 958
$ws.Range("D196").Value = 44232 ; $ws.Range("J196").Value = 3000 ; $ws.Range("K196").Value = 8500 ; $ws.Range("M196").Value = 8750 ; $ws.Range("P196").Value = 1458
$ws.Range("D197").Value = 44232 ; $ws.Range("J197").Value = 1600 ; $ws.Range("K197").Value = 7500 ; $ws.Range("L197").Value = 8000 ; $ws.Range("M197").Value = 7750 ; $ws.Range("P197").Value = 1292
$ws.Range("D198").Value = 44398 ; $ws.Range("J198").Value = 3400 ; $ws.Range("K198").Value = 8000 ; $ws.Range("L198").Value = 9000 ; $ws.Range("M198").Value = 8500 ; $ws.Range("P198").Value = 1417
$ws.Range("D199").Value = 44398 ; $ws.Range("J199").Value = 1640 ; $ws.Range("K199").Value = 6000 ; $ws.Range("L199").Value = 7000 ; $ws.Range("M199").Value = 6500 ; $ws.Range("P199").Value = 1083
$ws.Range("D200").Value = 44270 ; $ws.Range("J200").Value = 2600
$ws.Range("D201").Value = 44270 ; $ws.Range("J201").Value = 1400
$ws.Range("D202").Value = 44286 ; $ws.Range("K202").Value = 7000 ; $ws.Range("M202").Value = 7500 ; $ws.Range("P202").Value = 1250
$ws.Range("D203").Value = 44286 ; $ws.Range("J203").Value = 1660 ; $ws.Range("K203").Value = 5000 ; $ws.Range("L203").Value = 6000 ; $ws.Range("M203").Value = 5500 ; $ws.Range("P203").Value = 917
$ws.Range("D204").Value = 44244 ; $ws.Range("J204").Value = 3200 ; $ws.Range("K204").Value = 7500 ; $ws.Range("L204").Value = 8000 ; $ws.Range("M204").Value = 7750 ; $ws.Range("P204").Value = 1292
$ws.Range("D205").Value = 44244 ; $ws.Range("K205").Value = 6500 ; $ws.Range("L205").Value = 7000 ; $ws.Range("M205").Value = 6750 ; $ws.Range("P205").Value = 1125
$ws.Range("D206").Value = 44223 ; $ws.Range("J206").Value = 2900 ; $ws.Range("K206").Value = 9000 ; $ws.Range("L206").Value = 10000 ; $ws.Range("M206").Value = 9500 ; $ws.Range("P206").Value = 1583
$ws.Range("D207").Value = 44223 ; $ws.Range("J207").Value = 1600
$ws.Range("D208").Value = 44179 ; $ws.Range("J208").Value = 2600 ; $ws.Range("K208").Value = 8500 ; $ws.Range("L208").Value = 9000 ; $ws.Range("M208").Value = 8750 ; $ws.Range("P208").Value = 1458
$ws.Range("D209").Value = 44179 ; $ws.Range("J209").Value = 1540 ; $ws.Range("K209").Value = 7500 ; $ws.Range("L209").Value = 8000 ; $ws.Range("M209").Value = 7750 ; $ws.Range("P209").Value = 1292

# Append two brand-new rows (210, 211) holding the data that used to occupy rows 208-209.
$ws.Range("A210").Value = 8
$ws.Range("B210").Value = "Terminal La Palmera de La Serena"
$ws.Range("C210").Value = "Coquimbo"
$ws.Range("D210").Value = 44284
$ws.Range("D210").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E210").Value = 4
$ws.Range("F210").Value = 100112017
$ws.Range("G210").Value = "Apio"
$ws.Range("H210").Value = "Americana (o)"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 2800
$ws.Range("K210").Value = 7000
$ws.Range("L210").Value = 8000
$ws.Range("M210").Value = 7500
$ws.Range("N210").Value = "$/docena de matas"
$ws.Range("O210").Value = "Provincia del Elquí"
$ws.Range("P210").Value = 1250
$ws.Range("Q210").Value = 6
$ws.Range("R210").Value = "Hortaliza"

$ws.Range("A211").Value = 8
$ws.Range("B211").Value = "Terminal La Palmera de La Serena"
$ws.Range("C211").Value = "Coquimbo"
$ws.Range("D211").Value = 44284
$ws.Range("D211").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E211").Value = 4
$ws.Range("F211").Value = 100112017
$ws.Range("G211").Value = "Apio"
$ws.Range("H211").Value = "Americana (o)"
$ws.Range("I211").Value = "Segunda"
$ws.Range("J211").Value = 1460
$ws.Range("K211").Value = 5000
$ws.Range("L211").Value = 6000
$ws.Range("M211").Value = 5500
$ws.Range("N211").Value = "$/docena de matas"
$ws.Range("O211").Value = "Provincia del Elquí"
$ws.Range("P211").Value = 917
$ws.Range("Q211").Value = 6
$ws.Range("R211").Value = "Hortaliza"
